$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking price cells so exact string formatting
# (trailing zeros, decimal grouping style) is preserved instead of Excel
# auto-converting the entry to a Number.
$textCells = @("D5","D6","D7","D9","D10","D11","D12","D15","D19","D20","D22","D24","D26","D27","D28","D29","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D43","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values from the source diff.
$ws.Range("D2").Value = "26.284.96"
$ws.Range("D3").Value = "1.663.07"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("E4").Value = "  +0.72%  "
$ws.Range("D5").Value = "218.54"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").Value = "0.5317"
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("D7").Value = "1.010"
$ws.Range("E7").Value = "  +0.72%  "
$ws.Range("E8").Value = "  +0.94%  "
$ws.Range("D9").Value = "0.06364"
$ws.Range("E9").Value = "  +0.44%  "
$ws.Range("D10").Value = "20.52"
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("D11").Value = "0.07843"
$ws.Range("E11").Value = "  +1.15%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "4.553"
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.685.91"
$ws.Range("E13").Value = "  +1.81%  "
$ws.Range("D14").Value = "1.892.64"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").Value = "0.5536"
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("D16").Value = "0.0₅8185"
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("E18").Value = "  +0.69%  "
$ws.Range("D19").Value = "4.660"
$ws.Range("E19").Value = "  +2.62%  "
$ws.Range("D20").Value = "192.49"
$ws.Range("E20").Value = "  -0.72%  "
$ws.Range("E21").Value = "  +1.64%  "
$ws.Range("D22").Value = "6.055"
$ws.Range("E22").Value = "  +0.87%  "
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("D24").Value = "145.13"
$ws.Range("E24").Value = "  +3.29%  "
$ws.Range("E25").Value = "  -1.75%  "
$ws.Range("D26").Value = "7.235"
$ws.Range("E26").Value = "  -0.62%  "
$ws.Range("D27").Value = "16.13"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("D28").Value = "1.489"
$ws.Range("E28").Value = "  +3.95%  "
$ws.Range("D29").Value = "0.05877"
$ws.Range("E29").Value = "  -1.18%  "
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("D31").Value = "3.586"
$ws.Range("E31").Value = "  +2.02%  "
$ws.Range("D32").Value = "3.302"
$ws.Range("E32").Value = "  +2.05%  "
$ws.Range("D33").Value = "1.612"
$ws.Range("E33").Value = "  +4.20%  "
$ws.Range("D34").Value = "0.9591"
$ws.Range("E34").Value = "  +1.33%  "
$ws.Range("D35").Value = "2.817"
$ws.Range("E35").Value = "  +2.06%  "
$ws.Range("D36").Value = "2.425"
$ws.Range("E36").Value = "  +0.57%  "
$ws.Range("D37").Value = "0.5816"
$ws.Range("E37").Value = "  +3.28%  "
$ws.Range("D38").Value = "0.01613"
$ws.Range("E38").Value = "  +0.25%  "
$ws.Range("D39").Value = "5.897"
$ws.Range("E39").Value = "  +0.87%  "
$ws.Range("D40").Value = "0.8547"
$ws.Range("E40").Value = "  +0.96%  "
$ws.Range("E41").Value = "  +0.70%  "
$ws.Range("D42").Value = "1.047.46"
$ws.Range("E42").Value = "  +3.54%  "
$ws.Range("D43").Value = "104.24"
$ws.Range("E43").Value = "  +2.93%  "
$ws.Range("D44").Value = "1.805.35"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("E45").Value = "  +0.58%  "
$ws.Range("E46").Value = "  +2.08%  "
$ws.Range("D47").Value = "1.013"
$ws.Range("E47").Value = "  +1.03%  "
$ws.Range("D48").Value = "0.4374"
$ws.Range("E48").Value = "  +1.98%  "
$ws.Range("D49").Value = "7.953"
$ws.Range("E49").Value = "  +2.79%  "
$ws.Range("D50").Value = "0.05165"
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("D51").Value = "1.443"
$ws.Range("E51").Value = "  -1.67%  "
